$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mixer2-BOM")

# Row 5: merge "Bought" (D5=10) into "Have" (C5), was 2 -> now 4
$ws.Range("C5").Value = 4
$ws.Range("D5").ClearContents()

# Row 15: merge "Bought" (D15=20) into "Have" (C15), was 0 -> now 2
$ws.Range("C15").Value = 2
$ws.Range("D15").ClearContents()

# Row 21: merge "Bought" (D21=10) into "Have" (C21), was 0 -> now 1
$ws.Range("C21").Value = 1
$ws.Range("D21").ClearContents()

# Update selection to reflect last-edited cell
$ws.Range("D21").Select()
